# Update demo PPI workbook to be consistent with "main":
#  1. Change the corner header cell on the "network" sheet.
#  2. Add a new "optimization_parameters" sheet (right after "network")
#     describing the species/taxon/workbook-type metadata.

$wb = $excel.ActiveWorkbook

# --- 1. Update the corner/header cell on the existing "network" sheet ---
$ws1 = $wb.Worksheets.Item("network")
$ws1.Range("A1").Value = "cols protein1/ rows protein2"

# --- 2. Add the new "optimization_parameters" sheet right after "network" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "optimization_parameters"

$ws2.Range("A1").Value = "optimization_parameter"
$ws2.Range("B1").Value = "value"

$ws2.Range("A2").Value = "species"
$ws2.Range("B2").Value = "Saccharomyces cerevisiae"

$ws2.Range("A3").Value = "taxon_id"
$ws2.Range("B3").Value = 559292

$ws2.Range("A4").Value = "workbookType"
$ws2.Range("B4").Value = "protein-protein-physical-interaction"
